$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 text updates ---
$ws.Range("B2").Value = "hsdfh"
$ws.Range("C2").Value = "qhgfawgf"
$ws.Range("D2").Value = "hgfwhf"
$ws.Range("E2").Value = "ftqsf "
$ws.Range("F2").Value = "qfw"
$ws.Range("G2").Value = "q"
$ws.Range("H2").Value = "yf"
$ws.Range("I2").Value = "uqe"
$ws.Range("J2").Value = "utqe"

# --- Row 3 text updates ---
$ws.Range("B3").Value = "ghfaghfdyq"
$ws.Range("C3").Value = "qyit"
$ws.Range("D3").Value = "qdtiq"
$ws.Range("E3").Value = "it"
$ws.Range("F3").Value = "qe"
$ws.Range("G3").Value = "qdeyqedqoqe"
$ws.Range("H3").Value = "to"
$ws.Range("I3").Value = "qe"
$ws.Range("J3").Value = "qteq"

# --- Row 4 - new row of data ---
# A4 holds the text "3" (like A2/A3 which hold text "1"/"2"), so force
# text formatting before assigning a numeric-looking string.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "3"
$ws.Range("B4").Value = "wdqwdqwqw"
$ws.Range("C4").Value = "wereewt"
$ws.Range("D4").Value = "bashgdgf"
$ws.Range("E4").Value = "fgqwgf"
$ws.Range("F4").Value = "fg"
$ws.Range("G4").Value = "fg"
$ws.Range("H4").Value = "qgf"
$ws.Range("I4").Value = "qfg"
$ws.Range("J4").Value = "qfge"

# --- Selection moved back to A1 (e.g. after search bar usage) ---
$ws.Range("A1").Select()
